# Updated BOM with Digikey part numbers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quantity correction -------------------------------------------------
# Row 3 (0.22uF capacitor) quantity corrected from 1 to 3
$ws.Range("A3").Value = 3

# --- Part number correction ----------------------------------------------
# Row 34 (MIC803) device part number corrected
$ws.Range("C34").Value = "MIC803-31D2VM3TR"

# --- New Digikey part numbers in column G ---------------------------------
$ws.Range("G3").Value  = "587-1287-1-ND"
$ws.Range("G4").Value  = "399-1158-1-ND"
$ws.Range("G7").Value  = "399-1151-1-ND"
$ws.Range("G8").Value  = "399-7158-1-ND"
$ws.Range("G10").Value = "399-5620-1-ND"
$ws.Range("G14").Value = "408-1556-1-ND"
$ws.Range("G15").Value = "311-.75SCT-ND"
$ws.Range("G16").Value = "1276-6000-1-ND"
$ws.Range("G18").Value = "311-10KARCT-ND"
$ws.Range("G22").Value = "311-33.0KCRCT-ND"
$ws.Range("G24").Value = "311-4.70KCRCT-ND"
$ws.Range("G25").Value = "RHM470CHCT-ND"
$ws.Range("G27").Value = "311-68.0KCRCT-ND"
$ws.Range("G28").Value = "RHC2512FT7R50CT-ND"
$ws.Range("G30").Value = "ZXMP2120FFCT-ND"
$ws.Range("G34").Value = "576-3814-1-ND"
$ws.Range("G36").Value = "342-1082-1-ND"

# Some of the pasted-in Digikey values came in with Arial formatting
$arialCells = @("G8","G15","G16","G22","G25","G27","G28","G30","G36")
foreach ($addr in $arialCells) {
    $rng = $ws.Range($addr)
    $rng.Font.Name = "Arial"
    $rng.Font.Color = 0
}

# Row 35 (AB0805) Newark part number - "Newark " plain + bold "69W6470"
$ws.Range("G35").Value = "Newark 69W6470"
$ws.Range("G35").Characters(8, 7).Font.Bold = $true

# --- Restore last-used selection -----------------------------------------
$ws.Rows(42).Select()

echo "done"
